{"js": "// Add a new row to the (only) table in the document. The new cell holds a\n// paragraph built out of many runs -- one word/space at a time -- with the\n// word \"nonbreaking\" flanked by actual non-breaking space characters\n// (U+00A0) rather than regular spaces, matching the new test fixture added\n// for \"nonbreaking space removal\".\n\nconst table = context.document.body.tables.getFirst();\n\n// Add an empty row at the end; it automatically inherits the surrounding\n// rows' cell formatting (borders, trHeight) and paragraph formatting\n// (NormalWeb style, 11pt/sz=22).\ntable.addRows(\"End\", 1, [[\"\"]]);\nawait context.sync();\n\nconst row = table.rows.getLast();\nconst cell = row.cells.getFirst();\nconst para = cell.body.paragraphs.getFirst();\n\nconst NBSP = \"\\u00A0\";\n\n// Each word and each separating space is inserted as its own piece of text,\n// except the space before/after \"nonbreaking\", which are non-breaking\n// spaces (the whole point of this fixture).\nconst parts = [\n  \"This\",\n  \" \",\n  \"is\",\n  \" \",\n  \"a\",\n  \" \",\n  \"single\",\n  \" \",\n  \"run\",\n  \" \",\n  \"of\",\n  \" \",\n  \"text\",\n  \" \",\n  \"that\",\n  \" \",\n  \"only\",\n  \" \",\n  \"contains\",\n  NBSP + \"nonbreaking\" + NBSP,\n  \"spaces.\",\n];\n\nfor (const part of parts) {\n  para.insertText(part, \"End\");\n}\nawait context.sync();\n", "ps1": "# Add a new row to the (only) table in the document. The new cell holds a\n# paragraph built out of many runs -- one word/space at a time -- with the\n# word \"nonbreaking\" flanked by actual non-breaking space characters\n# (U+00A0) rather than regular spaces, matching the new test fixture added\n# for \"nonbreaking space removal\".\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Add an empty row at the end; it automatically inherits the surrounding\n# rows' cell formatting (borders, trHeight) and paragraph formatting\n# (NormalWeb style, 11pt/sz=22).\n$table.Rows.Add() | Out-Null\n\n$row = $table.Rows.Last\n$cell = $row.Cells.Item(1)\n$rng = $cell.Range\n$rng.Collapse(0) | Out-Null\n\n$nbsp = [char]0x00A0\n\n# Each word and each separating space is inserted as its own piece of text,\n# except the space before/after \"nonbreaking\", which are non-breaking\n# spaces (the whole point of this fixture).\n$parts = @(\n    \"This\",\n    \" \",\n    \"is\",\n    \" \",\n    \"a\",\n    \" \",\n    \"single\",\n    \" \",\n    \"run\",\n    \" \",\n    \"of\",\n    \" \",\n    \"text\",\n    \" \",\n    \"that\",\n    \" \",\n    \"only\",\n    \" \",\n    \"contains\",\n    \"$nbsp\" + \"nonbreaking\" + \"$nbsp\",\n    \"spaces.\"\n)\n\nforeach ($part in $parts) {\n    $rng.InsertAfter($part)\n    $rng.Collapse(0) | Out-Null\n}\n"}
